$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D28").Value = "[Pixhawk] Simulation Using Gazebo_ROS and Mavros"
$ws.Range("E28").Value = "https://ropiens.tistory.com/149"

$ws.Range("D29").Value = "[만화] 인턴일기 51~57"
$ws.Range("E29").Value = "https://blog.promedius.ai/intern-life-8/"

$ws.Range("D52").Value = "사건까지 걸린 시간은?"
